$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 6000
$ws.Range("B2").Value = 6000
$ws.Range("C2").Value = 10

$ws.Range("A3").Value = 5000
$ws.Range("B3").Value = 6000
$ws.Range("C3").Value = 10

$ws.Range("A4").Value = 5000
$ws.Range("B4").Value = 2000
$ws.Range("C4").Value = 10

$ws.Range("A5").Value = "'5000.0"
$ws.Range("B5").Value = "'7000.0"
$ws.Range("C5").Value = "'10"
